$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "95-6="
$t.Cell(1,2).Range.Text = "24-6="
$t.Cell(1,3).Range.Text = "71-3="
$t.Cell(1,4).Range.Text = "24+48="
$t.Cell(1,5).Range.Text = "62-48="
$t.Cell(2,1).Range.Text = "36+36="
$t.Cell(2,2).Range.Text = "39+34="
$t.Cell(2,3).Range.Text = "16+59="
$t.Cell(2,4).Range.Text = "98-49="
$t.Cell(2,5).Range.Text = "83-4="
$t.Cell(3,1).Range.Text = "46-18="
$t.Cell(3,2).Range.Text = "58+18="
$t.Cell(3,3).Range.Text = "90-31="
$t.Cell(3,4).Range.Text = "83-36="
$t.Cell(3,5).Range.Text = "81-22="
$t.Cell(4,1).Range.Text = "47+39="
$t.Cell(4,2).Range.Text = "8+48="
$t.Cell(4,3).Range.Text = "63+28="
$t.Cell(4,4).Range.Text = "24+7="
$t.Cell(4,5).Range.Text = "29+54="
$t.Cell(5,1).Range.Text = "71-39="
$t.Cell(5,2).Range.Text = "64-27="
$t.Cell(5,3).Range.Text = "16+79="
$t.Cell(5,4).Range.Text = "60-5="
$t.Cell(5,5).Range.Text = "92-9="
$t.Cell(6,1).Range.Text = "73-67="
$t.Cell(6,2).Range.Text = "17+17="
$t.Cell(6,3).Range.Text = "87-49="
$t.Cell(6,4).Range.Text = "59+28="
$t.Cell(6,5).Range.Text = "83-25="
$t.Cell(7,1).Range.Text = "93-89="
$t.Cell(7,2).Range.Text = "81-22="
$t.Cell(7,3).Range.Text = "93-45="
$t.Cell(7,4).Range.Text = "11-8="
$t.Cell(7,5).Range.Text = "33+29="
$t.Cell(8,1).Range.Text = "62-49="
$t.Cell(8,2).Range.Text = "64-35="
$t.Cell(8,3).Range.Text = "45+6="
$t.Cell(8,4).Range.Text = "8+17="
$t.Cell(8,5).Range.Text = "36+45="
$t.Cell(9,1).Range.Text = "42-37="
$t.Cell(9,2).Range.Text = "55-16="
$t.Cell(9,3).Range.Text = "27+49="
$t.Cell(9,4).Range.Text = "52-16="
$t.Cell(9,5).Range.Text = "23+39="
$t.Cell(10,1).Range.Text = "8+85="
$t.Cell(10,2).Range.Text = "61-15="
$t.Cell(10,3).Range.Text = "56+38="
$t.Cell(10,4).Range.Text = "81-57="
$t.Cell(10,5).Range.Text = "38+29="
$t.Cell(11,1).Range.Text = "61-43="
$t.Cell(11,2).Range.Text = "50-3="
$t.Cell(11,3).Range.Text = "55+9="
$t.Cell(11,4).Range.Text = "69+27="
$t.Cell(11,5).Range.Text = "42-4="
$t.Cell(12,1).Range.Text = "94-76="
$t.Cell(12,2).Range.Text = "69+22="
$t.Cell(12,3).Range.Text = "22-19="
$t.Cell(12,4).Range.Text = "9+89="
$t.Cell(12,5).Range.Text = "90-13="
$t.Cell(13,1).Range.Text = "83+8="
$t.Cell(13,2).Range.Text = "96-39="
$t.Cell(13,3).Range.Text = "55+37="
$t.Cell(13,4).Range.Text = "39+7="
$t.Cell(13,5).Range.Text = "37+34="
$t.Cell(14,1).Range.Text = "6+28="
$t.Cell(14,2).Range.Text = "57-38="
$t.Cell(14,3).Range.Text = "62-27="
$t.Cell(14,4).Range.Text = "10-8="
$t.Cell(14,5).Range.Text = "62-54="
$t.Cell(15,1).Range.Text = "63+8="
$t.Cell(15,2).Range.Text = "61-13="
$t.Cell(15,3).Range.Text = "70-62="
$t.Cell(15,4).Range.Text = "58+3="
$t.Cell(15,5).Range.Text = "27-8="
$t.Cell(16,1).Range.Text = "91-19="
$t.Cell(16,2).Range.Text = "71-34="
$t.Cell(16,3).Range.Text = "8+25="
$t.Cell(16,4).Range.Text = "14+78="
$t.Cell(16,5).Range.Text = "7+16="
$t.Cell(17,1).Range.Text = "85-67="
$t.Cell(17,2).Range.Text = "27+4="
$t.Cell(17,3).Range.Text = "45-38="
$t.Cell(17,4).Range.Text = "17+74="
$t.Cell(17,5).Range.Text = "19+5="
$t.Cell(18,1).Range.Text = "67-19="
$t.Cell(18,2).Range.Text = "91-74="
$t.Cell(18,3).Range.Text = "51-15="
$t.Cell(18,4).Range.Text = "3+59="
$t.Cell(18,5).Range.Text = "90-65="
$t.Cell(19,1).Range.Text = "46+46="
$t.Cell(19,2).Range.Text = "44-18="
$t.Cell(19,3).Range.Text = "19+75="
$t.Cell(19,4).Range.Text = "37+34="
$t.Cell(19,5).Range.Text = "7+27="
$t.Cell(20,1).Range.Text = "68-49="
$t.Cell(20,2).Range.Text = "61-18="
$t.Cell(20,3).Range.Text = "13+68="
$t.Cell(20,4).Range.Text = "81-5="
$t.Cell(20,5).Range.Text = "7+54="
